# Updates cryptos list values (price & 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.468.86"
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = "'3.080.80"
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'545.52"
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("D6").Value = "'139.11"
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = "'3.075.39"
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").Value = "'0.158"
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("D11").Value = "'6.43"
$ws.Range("E11").Value = '  +2.27%  '
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = '  -2.94%  '
$ws.Range("E13").Value = '  +3.41%  '
$ws.Range("D14").Value = "'35.01"
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("D15").Value = "'3.579.06"
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = "'63.437.35"
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").Value = "'3.081.11"
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("E19").Value = '  -1.26%  '
$ws.Range("D20").Value = "'476.72"
$ws.Range("E20").Value = '  -2.57%  '
$ws.Range("D21").Value = "'13.52"
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("D22").Value = "'0.702"
$ws.Range("E22").Value = '  -2.42%  '
$ws.Range("D23").Value = "'7.09"
$ws.Range("E23").Value = '  -2.63%  '
$ws.Range("D24").Value = "'78.72"
$ws.Range("E24").Value = '  -0.59%  '
$ws.Range("D25").Value = "'12.23"
$ws.Range("E25").Value = '  -1.22%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("E28").Value = '  -6.84%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = "'26.28"
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("E31").Value = '  +3.33%  '
$ws.Range("E32").Value = '  -3.83%  '
$ws.Range("D33").Value = "'59.18"
$ws.Range("E33").Value = '  +0.91%  '
$ws.Range("D34").Value = "'2.32"
$ws.Range("E34").Value = '  -7.43%  '
$ws.Range("E35").Value = '  +7.46%  '
$ws.Range("D36").Value = "'6.03"
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("D37").Value = "'489.23"
$ws.Range("E37").Value = '  -3.88%  '
$ws.Range("D38").Value = "'3.278.51"
$ws.Range("E38").Value = '  +4.01%  '
$ws.Range("D39").Value = "'0.0403"
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("D40").Value = "'0.0798"
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("E41").Value = '  -1.24%  '
$ws.Range("D42").Value = "'8.18"
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("E44").Value = '  -2.04%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = "'25.44"
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("D47").Value = "'124.32"
$ws.Range("E47").Value = '  +2.95%  '
$ws.Range("E48").Value = '  -2.06%  '
$ws.Range("D49").Value = "'0.0₃0531"
$ws.Range("E49").Value = '  +4.58%  '
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("E51").Value = '  -0.46%  '
